$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data held in row 2 and row 3 for the changed columns
# (Date, Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)

$ws.Range("D2").Value2 = 44875
$ws.Range("J2").Value2 = 1000
$ws.Range("K2").Value2 = 1600
$ws.Range("L2").Value2 = 1700
$ws.Range("M2").Value2 = 1650
$ws.Range("P2").Value2 = 1650

$ws.Range("D3").Value2 = 44547
$ws.Range("J3").Value2 = 400
$ws.Range("K3").Value2 = 1500
$ws.Range("L3").Value2 = 1600
$ws.Range("M3").Value2 = 1550
$ws.Range("P3").Value2 = 1550
